# Add new columns I (I0) and J (IF) with per-row values, mirroring H's "1" style
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells, matching the other header cells' look (bold, centered, top-aligned, thin border)
$headerCells = @($ws.Cells.Item(1, 9), $ws.Cells.Item(1, 10))

$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

foreach ($cell in $headerCells) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108  # xlCenter
    $cell.VerticalAlignment = -4160    # xlTop
    $cell.Borders.LineStyle = 1        # xlContinuous (thin)
}

# Data values for rows 2..66 (column I == column J for every row)
$values = @(8, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 10, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 8, 9, 9, 9, 9, 9, 9, 9, 9, 10, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 10, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 5, 4, 4, 5)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i]
    $ws.Cells.Item($row, 10).Value = $values[$i]
}
